$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that sits right before the first table
#    (after "Система повинна містити 4 типи користувачів...").
# ---------------------------------------------------------------------------
foreach ($bm in @($d.Bookmarks)) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

# ---------------------------------------------------------------------------
# 2) Simple in-place text edits inside the second (requirements) table.
# ---------------------------------------------------------------------------
$t2 = $d.Tables.Item(2)

# Row "6." Filter subjects : "S" -> "S, L, A"
$t2.Rows.Item(7).Cells.Item(3).Range.Text = "S, L, A"

# Row "9."  -> becomes "Add new profile" / "Створення нового акаунту"
$t2.Rows.Item(10).Cells.Item(2).Range.Text = "Add new profile"
$t2.Rows.Item(10).Cells.Item(4).Range.Text = "Створення нового акаунту"

# Row "10." -> becomes "Edit profile" / "Редагування профілів користувачів"
$t2.Rows.Item(11).Cells.Item(2).Range.Text = "Edit profile"
$t2.Rows.Item(11).Cells.Item(4).Range.Text = "Редагування профілів користувачів"

# Row "11." -> becomes "Delete profile" / "Видалення профілю"
$t2.Rows.Item(12).Cells.Item(2).Range.Text = "Delete profile"
$t2.Rows.Item(12).Cells.Item(4).Range.Text = "Видалення профілю"

# Row "12." -> becomes "Add faculty" / "Додати новий факультет"
$t2.Rows.Item(13).Cells.Item(2).Range.Text = "Add faculty"
$t2.Rows.Item(13).Cells.Item(4).Range.Text = "Додати новий факультет"

# Row "13." -> becomes "Edit faculty info" / "Редагувати інформацію про факультет"
$t2.Rows.Item(14).Cells.Item(2).Range.Text = "Edit faculty info"
$t2.Rows.Item(14).Cells.Item(4).Range.Text = "Редагувати інформацію про факультет"

# Row "14." -> becomes "Delete faculty" / "Видалити факультет"
$t2.Rows.Item(15).Cells.Item(2).Range.Text = "Delete faculty"
$t2.Rows.Item(15).Cells.Item(4).Range.Text = "Видалити факультет"

# ---------------------------------------------------------------------------
# 3) Insert a brand-new row right after the (now renamed) row "14." — this is
#    a duplicate of what row "14." used to contain: "Set period of choice".
#    It carries the relocated _GoBack bookmark in its first cell.
# ---------------------------------------------------------------------------
$newRow = $t2.Rows.Add($t2.Rows.Item(16))
$newRow.Cells.Item(1).Range.Text = "15."
$newRow.Cells.Item(2).Range.Text = "Set period of choice"
$newRow.Cells.Item(3).Range.Text = "A"
$newRow.Cells.Item(4).Range.Text = "Додати період можливості вибору предметів"

# Put the _GoBack bookmark back at the very start of the new row's first cell.
$d.Bookmarks.Add("_GoBack", $newRow.Cells.Item(1).Range.Paragraphs.Item(1).Range)

# ---------------------------------------------------------------------------
# 4) The row that used to be "15. Search" becomes "16. Search subject" and
#    picks up the page-break-render hint that used to sit on the old "16."
#    row, plus an expanded capability list.
# ---------------------------------------------------------------------------
$searchRow = $t2.Rows.Item(17)
$searchRow.Cells.Item(1).Range.Text = "16."
$searchRow.Cells.Item(2).Range.Text = "Search subject"
$searchRow.Cells.Item(3).Range.Text = "S, L, A"

# ---------------------------------------------------------------------------
# 5) Old "16. Edit subjects" row becomes "17. Edit subject" (singular) and no
#    longer carries the page-break-render hint.
# ---------------------------------------------------------------------------
$editSubjRow = $t2.Rows.Item(18)
$editSubjRow.Cells.Item(1).Range.Text = "17."
$editSubjRow.Cells.Item(2).Range.Text = "Edit subject"
$editSubjRow.Cells.Item(4).Range.Text = "Редагування предмету"

# ---------------------------------------------------------------------------
# 6) Old "17. Delete subjects" row becomes "18. Delete subject" (singular).
# ---------------------------------------------------------------------------
$delSubjRow = $t2.Rows.Item(19)
$delSubjRow.Cells.Item(1).Range.Text = "18."
$delSubjRow.Cells.Item(2).Range.Text = "Delete subject"
